{"js": "// The CV's Education entry \"M.A. in Cognitive Psychology (Ph.D. expected\n// early 2022)\" originally lists a \"May 2019\" start date in its right-hand\n// date column. The edit changes that date to \"  Dec 2019\" (the two leading\n// spaces are part of the original diff's replacement text).\n//\n// Anchor on the unique phrase \"Ph.D. expected\" to find the correct\n// paragraph (there is another, unrelated \"May 2019\" elsewhere in the\n// document), then search for \"May 20\" inside that paragraph only and\n// replace it (together with the following \"19\") with \"  Dec 2019\".\n\nconst anchorResults = context.document.body.search(\"Ph.D. expected\", { matchCase: true });\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length === 0) {\n  throw new Error('Anchor text \"Ph.D. expected\" not found.');\n}\n\nconst targetParagraph = anchorResults.items[0].paragraphs.getFirst();\n\nconst dateResults = targetParagraph.search(\"May 2019\", { matchCase: true });\ndateResults.load(\"items,text\");\nawait context.sync();\n\nif (dateResults.items.length === 0) {\n  throw new Error('\"May 2019\" not found in the target paragraph.');\n}\n\nconst dateRange = dateResults.items[0];\ndateRange.insertText(\"  Dec 2019\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The CV's Education entry \"M.A. in Cognitive Psychology (Ph.D. expected\n# early 2022)\" originally lists a \"May 2019\" start date in its right-hand\n# date column. The edit changes that date to \"  Dec 2019\" (the two leading\n# spaces are part of the original replacement text).\n#\n# There is a second, unrelated \"May 2019\" elsewhere in the document (in the\n# \"Molly McKinney ... Sep. 2018 - May 2019\" line), so we must not do a\n# blind document-wide replace. Instead we anchor on the unique phrase\n# \"Ph.D. expected\" to locate the correct paragraph, then restrict the\n# Find/Replace to that paragraph's Range only.\n\n$d = $word.ActiveDocument\n\n$anchor = $d.Content\n$anchorFind = $anchor.Find\n$anchorFind.Text = \"Ph.D. expected\"\n$anchorFind.MatchCase = $true\n$found = $anchorFind.Execute()\n\nif (-not $found) {\n    throw 'Anchor text \"Ph.D. expected\" not found.'\n}\n\n$targetParagraph = $anchor.Paragraphs(1).Range\n\n$find = $targetParagraph.Find\n$find.Text = \"May 2019\"\n$find.Replacement.Text = \"  Dec 2019\"\n$find.MatchCase = $true\n$find.Forward = $true\n$find.Wrap = 0      # wdFindStop - do not search beyond the paragraph range\n$find.Format = $false\n\n$replaced = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 1)\n\nif (-not $replaced) {\n    throw '\"May 2019\" not found in the target paragraph.'\n}\n"}
